# Update the workbook for the May 2018 (2 Periods) weekend-weekday validation run.
$wb = $excel.ActiveWorkbook

# 1. Table of Contents: fix the description of "initstorage" (row 27, column E)
#    from "1st March 2018" to "1st May 2018".
$toc = $wb.Worksheets.Item("Table of Contents")
$toc.Range("E27").Value = "Initial reservoir storage on 1st May 2018 (acre-ft)"

# 2. "evap" sheet: the evaporation series (B4:B34) is replaced with the values
#    that used to live on the "Inflow" sheet.
$evapValues = @(931,930,930,930,931,931,931,931,931,931,931,931,931,932,933,933,934,935,937,938,939,939,939,940,940,941,942,943,944,945,946)
$evapSheet = $wb.Worksheets.Item("evap")
for ($i = 0; $i -lt $evapValues.Length; $i++) {
    $row = 4 + $i
    $evapSheet.Cells.Item($row, 2).Value = $evapValues[$i]
}

# 3. "Inflow" sheet: replace the inflow series (B4:B34) with the new recorded
#    inflow values for May 2018.
$inflowValues = @(7463,9189,14032,12362,13459,13430,13958,13490,10179,12774,10160,9988,14362,17827,19566,18601,19416,21007,20595,18950,19500,14353,18869,17428,15285,16056,17703,16804,23322,19487,18418)
$inflowSheet = $wb.Worksheets.Item("Inflow")
for ($i = 0; $i -lt $inflowValues.Length; $i++) {
    $row = 4 + $i
    $inflowSheet.Cells.Item($row, 2).Value = $inflowValues[$i]
}

# 4. "Scalar" sheet: update the derived/reported values that change as a
#    result of the new inflow/evaporation data (Sstore, Avail_Water, storage,
#    EQ1__ResMassBal, EQ2__reqpowerstorage, EQ3__maxstor, EQa_Inflow).
$scalar = $wb.Worksheets.Item("Scalar")
$scalar.Range("B12").Value = 12870148.45184   # Sstore
$scalar.Range("C22").Value = 13631127.736     # Avail_Water
$scalar.Range("C25").Value = 12870148.45184   # storage
$scalar.Range("B32").Value = -29000           # EQ1__ResMassBal
$scalar.Range("C32").Value = -29000
$scalar.Range("D32").Value = -29000
$scalar.Range("C33").Value = 12870148.45184   # EQ2__reqpowerstorage
$scalar.Range("C34").Value = 12870148.45184   # EQ3__maxstor
$scalar.Range("B36").Value = 13631127.736     # EQa_Inflow
$scalar.Range("C36").Value = 13631127.736
$scalar.Range("D36").Value = 13631127.736
